# "Generate Report for Archive"
#
# 1. The status text "Ready for handoff" becomes "In Translation" everywhere
#    it is used (Overview!E2/F2 mirror the per-locale Status cells
#    zh-cn!C2 and de-de!C2).
# 2. The Status-related columns shrink to fit the new, shorter text
#    (Overview columns E & F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = 12.5
